# Update gh-pages to output generated at 456a3b4
# Applies updated values to column F ("人气"/popularity-like counter) across
# the "展览", "演出" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
  5  = 177
  6  = 1055
  7  = 1036
  8  = 8044
  9  = 130
  10 = 194
  11 = 6846
  12 = 163
  13 = 299
  14 = 4911
  17 = 5334
  19 = 321
  20 = 323
  21 = 442
  22 = 310
  27 = 9038
  28 = 69
  29 = 1615
  33 = 833
  34 = 71
  36 = 1005
  37 = 1163
  38 = 53
  39 = 4706
  40 = 29
  41 = 374
  42 = 1156
  45 = 70
  46 = 32
  47 = 1237
  48 = 34
}
foreach ($row in $updates1.Keys) {
  $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$updates2 = @{
  3  = 38
  8  = 33
  17 = 886
}
foreach ($row in $updates2.Keys) {
  $ws2.Range("F$row").Value = $updates2[$row]
}

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
  6  = 177
  7  = 38
  8  = 1055
  9  = 1036
  10 = 8044
  11 = 130
  12 = 194
  13 = 6846
  14 = 163
  15 = 299
  17 = 4911
  19 = 5334
  20 = 1068
  21 = 321
  22 = 323
  23 = 442
  24 = 310
  30 = 9038
  31 = 69
  32 = 1615
  35 = 833
  36 = 71
  38 = 1005
  39 = 1163
  40 = 53
  41 = 4706
  42 = 374
  43 = 1156
  45 = 70
  46 = 32
  47 = 1237
  48 = 34
}
foreach ($row in $updates4.Keys) {
  $ws4.Range("F$row").Value = $updates4[$row]
}
